$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.169.12'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.840.90'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '361.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.02%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +4.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.06'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("E14").Value = '  +2.62%  '
$ws.Range("D15").Value = '3.286.83'
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("D16").Value = '2.836.97'
$ws.Range("E16").Value = '  +2.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.912'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.26%  '
$ws.Range("D18").Value = '52.115.74'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("E19").Value = '  +8.65%  '
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000100'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0488'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +29.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '53.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.29%  '
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("E34").Value = '  +2.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0846'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.53%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.34%  '
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '126.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.91%  '
$ws.Range("E45").Value = '  -3.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.37%  '
$ws.Range("D47").Value = '2.114.81'
$ws.Range("E47").Value = '  +1.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.991'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.14%  '
